$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '30.837.31'
$ws.Cells.Item(2, 5).Value = '  +0.67%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.890.05'
$ws.Cells.Item(3, 5).Value = '  +1.12%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.9996'
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  -0.03%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '240.06'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +2.05%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  -0.01%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.4790'
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  +1.72%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.2953'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +6.59%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.06630'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +3.82%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '18.85'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +4.83%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '101.40'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +18.86%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.890.10'
$ws.Cells.Item(12, 5).Value = '  +1.17%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.07644'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +2.49%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '5.143'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +3.04%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.6614'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +3.83%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '304.24'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +26.11%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '30.798.30'
$ws.Cells.Item(17, 5).Value = '  +0.67%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '13.19'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +2.31%  '

# Row 19
$ws.Cells.Item(19, 2).Value = 'ShibaInu'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '0.000007635'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +3.23%  '

# Row 20
$ws.Cells.Item(20, 2).Value = 'Dai'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '1.000'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +0.00%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '2.142.85'
$ws.Cells.Item(21, 5).Value = '  +2.22%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '1.002'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +0.08%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '5.184'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +3.64%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '6.211'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +2.45%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '169.14'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +2.02%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '9.356'
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -0.64%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '20.80'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +14.06%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '1.971'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +3.97%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.1123'
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +9.77%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '1.347'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -2.62%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '4.201'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +2.35%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '4.009'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +3.43%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.05094'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +3.19%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.7493'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +5.18%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '1.163'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +0.76%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.723'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +0.52%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.02004'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +4.84%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '2.707'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +0.42%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '2.060'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +3.13%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.8823'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -0.06%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'Quant'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '109.47'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +3.36%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +0.05%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.4218'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +2.44%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '5.661'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +1.69%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '68.31'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +8.48%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '7.389'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -0.28%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Algorand'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.1233'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -0.04%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '9.101'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +5.38%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '35.01'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +3.73%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.05655'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +1.46%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '1.401'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +1.39%  '
